$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12 -> subject 10
$ws.Range("B12").Value = "G"
$ws.Range("C12").Value = 8
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 11
$ws.Range("F12").Value = 12
$ws.Range("G12").Value = 12
$ws.Range("H12").Value = 3

# Row 13 -> subject 11
$ws.Range("B13").Value = "F"
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 4
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 6

# Update the active selection to match the diff
$ws.Range("I12").Select()
